$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 221, shifting the existing rows 221-238 down to 222-239.
$ws.Rows.Item(221).Insert()

# Populate the new row 221 with the new record's data.
$ws.Cells.Item(221, 1).Value = 1
$ws.Cells.Item(221, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(221, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(221, 4).Value = 45013
$ws.Cells.Item(221, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(221, 5).Value = 15
$ws.Cells.Item(221, 6).Value = "Fruta"
$ws.Cells.Item(221, 7).Value = 100108
$ws.Cells.Item(221, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(221, 9).Value = 100108002
$ws.Cells.Item(221, 10).Value = "Mango"
$ws.Cells.Item(221, 11).Value = "Sin especificar"
$ws.Cells.Item(221, 12).Value = "Segunda"
$ws.Cells.Item(221, 13).Value = 1000
$ws.Cells.Item(221, 14).Value = 4500
$ws.Cells.Item(221, 15).Value = 4800
$ws.Cells.Item(221, 16).Value = 4635
$ws.Cells.Item(221, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(221, 18).Value = "Perú"
$ws.Cells.Item(221, 19).Value = 1159
$ws.Cells.Item(221, 20).Value = 4
